$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("语义分析")
$ws2.Range("B2").Value = 9
$ws2.Activate()
$ws2.Range("C4").Select()
